# Applies the 05-11-2023 betexplorer scrape update:
#  - 5 pairs of rows had their match data (cols F:V) swapped back into the
#    correct chronological/pairing order (A:E - index/country/league/season/date -
#    stay put on each row).
#  - 3 brand-new match rows (93-95) are appended after the previous last row (92).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, [int]$rowA, [int]$rowB)

    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    foreach ($col in $cols) {
        $refA = "$col$rowA"
        $refB = "$col$rowB"
        $valA = $ws.Range($refA).Value2
        $valB = $ws.Range($refB).Value2
        $ws.Range($refA).Value = $valB
        $ws.Range($refB).Value = $valA
    }
}

# --- Swap the five mis-paired match rows back into place ---
Swap-RowData $ws 6 7
Swap-RowData $ws 44 45
Swap-RowData $ws 47 48
Swap-RowData $ws 50 51
Swap-RowData $ws 53 54

# --- Append the three newly scraped matches (rows 93-95) ---

function Add-MatchRow {
    param(
        $ws, [int]$row, [int]$idx, [string]$country, [string]$league, [string]$season,
        [double]$date, [string]$home, [int]$homeGoals, [string]$away, [int]$awayGoals,
        [double]$odd1Open, [string]$odd1OpenDate, [double]$odd1Close, [string]$odd1CloseDate,
        [double]$oddXOpen, [string]$oddXOpenDate, [double]$oddXClose, [string]$oddXCloseDate,
        [double]$odd2Open, [string]$odd2OpenDate, [double]$odd2Close, [string]$odd2CloseDate,
        [string]$url
    )

    # Column A: bold/centered/bordered index style (style id 1 in the sheet)
    $ws.Range("A92").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $idx

    $ws.Range("B$row").Value = $country
    $ws.Range("C$row").Value = $league
    $ws.Range("D$row").Value = $season

    # Column E: date/time number format (style id 2 in the sheet)
    $ws.Range("E92").Copy()
    $ws.Range("E$row").PasteSpecial(-4122)
    $ws.Range("E$row").Value = $date

    $ws.Range("F$row").Value = $home
    $ws.Range("G$row").Value = $homeGoals
    $ws.Range("H$row").Value = $away
    $ws.Range("I$row").Value = $awayGoals

    $ws.Range("J$row").Value = $odd1Open
    $ws.Range("K$row").Value = $odd1OpenDate
    $ws.Range("L$row").Value = $odd1Close
    $ws.Range("M$row").Value = $odd1CloseDate

    $ws.Range("N$row").Value = $oddXOpen
    $ws.Range("O$row").Value = $oddXOpenDate
    $ws.Range("P$row").Value = $oddXClose
    $ws.Range("Q$row").Value = $oddXCloseDate

    $ws.Range("R$row").Value = $odd2Open
    $ws.Range("S$row").Value = $odd2OpenDate
    $ws.Range("T$row").Value = $odd2Close
    $ws.Range("U$row").Value = $odd2CloseDate

    $ws.Range("V$row").Value = $url
}

Add-MatchRow $ws 93 92 "turkey" "1-lig" "2023-2024" 45234.47916666666 `
    "Corum" 2 "Keciorengucu" 3 `
    1.97 "30/10/2023 14:42" 2.08 "04/11/2023 10:56" `
    3.42 "30/10/2023 14:42" 3.39 "04/11/2023 10:56" `
    3.92 "30/10/2023 14:42" 3.7 "04/11/2023 10:56" `
    "https://www.betexplorer.com/football/turkey/1-lig/corum-fk-keciorengucu/A9vDdith/"

Add-MatchRow $ws 94 93 "turkey" "1-lig" "2023-2024" 45234.58333333334 `
    "Genclerbirligi" 0 "Boluspor" 0 `
    1.84 "28/10/2023 18:13" 1.97 "04/11/2023 13:56" `
    3.6 "28/10/2023 18:13" 3.27 "04/11/2023 13:56" `
    4.27 "28/10/2023 18:13" 4.29 "04/11/2023 13:56" `
    "https://www.betexplorer.com/football/turkey/1-lig/genclerbirligi-boluspor/bVAfvAQN/"

Add-MatchRow $ws 95 94 "turkey" "1-lig" "2023-2024" 45234.70833333334 `
    "Goztepe" 0 "Kocaelispor" 1 `
    2.02 "29/10/2023 14:42" 2.14 "04/11/2023 16:51" `
    3.3 "29/10/2023 14:42" 3.23 "04/11/2023 16:51" `
    3.88 "29/10/2023 14:42" 3.73 "04/11/2023 16:51" `
    "https://www.betexplorer.com/football/turkey/1-lig/goztepe-kocaelispor/0bXPgkQA/"
